$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuels")
Write-Host $ws.Name
